$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create row 22 by copying row 21 (values, formulas, formats) and inserting
# it right below, so styles (s="7" etc.), number formats and the existing
# "error" cells (G/H as #N/A) come along for the ride exactly as they look
# on row 21.
$ws.Range("A21:M21").Copy() | Out-Null
$ws.Range("A22:M22").Insert(-4121) | Out-Null
$ws.Application.CutCopyMode = 0

# Fix up the formulas that need to reference the new previous row.
$ws.Range("A22").Formula = "=A21+1"
$ws.Range("B22").Formula = "=B21+1"
$ws.Range("K22").Formula = "=K21+L22"

# New day's reported figures (2020-04-06).
$ws.Range("C22").Value2 = 430
$ws.Range("D22").Value2 = 102
$ws.Range("E22").Value2 = 1456
$ws.Range("F22").Value2 = 1988
$ws.Range("I22").Value2 = 0
$ws.Range("J22").Value2 = 3
$ws.Range("L22").Value2 = 0
$ws.Range("M22").Value2 = 66

# G22/H22 stay as the same "#N/A" placeholders that row 21 already had,
# copied over by the Insert above.

# Match the saved selection from the authored workbook.
$ws.Range("M23").Select() | Out-Null
